$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 21-23, mirroring row 4's layout but with the new save-location string
# and new "Sampling Rate Fluor" column (J)
$mice = @("R3M2368", "R3M2370", "R3M2372")
for ($i = 0; $i -lt $mice.Length; $i++) {
    $row = 21 + $i
    $ws.Cells.Item($row, 1).Value = 190530
    $ws.Cells.Item($row, 2).Value = $mice[$i]
    $ws.Cells.Item($row, 3).Value = "D:\RawData"
    $ws.Cells.Item($row, 4).Value = "D:\ProcessedData\FAD-GcAMP-New"
    $ws.Cells.Item($row, 5).Value = "EastOIS1_Fluor"
    $ws.Cells.Item($row, 6).Value = "{'fc'}"
    $ws.Cells.Item($row, 7).Value = 20
    $ws.Cells.Item($row, 8).Value = 5
    $ws.Cells.Item($row, 9).Value = 9
    $ws.Cells.Item($row, 10).Value = 9
}

# New column J header
$ws.Range("J1").Value = "Sampling Rate Fluor"

$ws.Columns.Item(10).AutoFit()

$ws.Range("B25").Select()
